$wb = $excel.ActiveWorkbook

# --- Sheet "Log" (sheet1): append rows 15-17 ---
$log = $wb.Worksheets.Item("Log")
$log.Range("A15").Value = '27/02/2025'
$log.Range("B15").Value = '15:00'
$log.Range("C15").Value = 'Exportar a Excel'
$log.Range("D15").Value = 'Botón "Exportar a Excel" con icono (mismo estilo que los del modal: gris, sencillo). Exporta la tabla de transacciones tal como está en Supabase: todas las columnas (fecha, mes, anio, tipo_movimiento, monto, status, medio_pago, descripcion, cliente, categoria, cat_desc, origen_archivo, cuenta_contable) en una hoja Excel para poder analizar los datos desde Excel. Librería SheetJS (xlsx) en el navegador.'
$log.Range("E15").Value = 'Diagnostico'
$log.Range("A16").Value = '27/02/2025'
$log.Range("B16").Value = '15:10'
$log.Range("C16").Value = 'Exportar transacciones crudas'
$log.Range("D16").Value = 'Ajuste: el botón Exportar a Excel pasa a exportar directamente la tabla de transacciones (datos crudos de Supabase), no el resumen flujo por mes, para permitir manipular y analizar los datos desde Excel.'
$log.Range("E16").Value = 'Diagnostico'
$log.Range("A17").Value = '27/02/2025'
$log.Range("B17").Value = '15:30'
$log.Range("C17").Value = 'Regla flujo despliegue y versiones'
$log.Range("D17").Value = 'Nueva regla: al final de cada tarea el usuario prueba en local y confirma; recién entonces el asistente despliega (git push). Se agrega hoja Versiones en la bitácora para registrar versión incremental en cada despliegue (1.0, 1.1, …).'
$log.Range("E17").Value = 'Diagnostico'

# --- Sheet "Resumen" (sheet2): append rows 22-24 ---
$resumen = $wb.Worksheets.Item("Resumen")
$resumen.Range("A22").Value = 'Exportar a Excel'
$resumen.Range("B22").Value = 'Botón en la barra de la tabla (solo icono). Exporta la tabla de transacciones tal como está en Supabase: una hoja "Transacciones" con columnas fecha, mes, anio, tipo_movimiento, monto, status, medio_pago, descripcion, cliente, categoria, cat_desc, origen_archivo, cuenta_contable. Permite analizar y manipular los datos desde Excel.'
$resumen.Range("A23").Value = 'Flujo de despliegue'
$resumen.Range("B23").Value = 'Al terminar cada tarea: el usuario prueba en local y confirma; recién entonces el asistente hace git add, commit y push (Vercel redepliega automático). No se despliega hasta confirmación.'
$resumen.Range("A24").Value = 'Versiones en bitácora'
$resumen.Range("B24").Value = 'Hoja "Versiones" en Bitacora_tareas.xlsx: registro incremental (1.0, 1.1, …) con fecha y descripción de cada despliegue a Git/Vercel.'

# --- New sheet "Versiones" appended after "Ref Git y Vercel" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$versiones = $wb.Worksheets.Add($null, $lastSheet)
$versiones.Name = "Versiones"

$versiones.Columns.Item(1).ColumnWidth = 7.998697916666667
$versiones.Columns.Item(2).ColumnWidth = 11.998697916666666
$versiones.Columns.Item(3).ColumnWidth = 74.99869791666667

$versiones.Range("A1").Value = 'Versión'
$versiones.Range("B1").Value = 'Fecha'
$versiones.Range("C1").Value = 'Descripción'
$versiones.Range("A2").NumberFormat = "@"
$versiones.Range("A2").Value = '1.0'
$versiones.Range("A2").Style = "Normal"
$versiones.Range("B2").Value = '27/02/2025'
$versiones.Range("C2").Value = 'Estado inicial: dashboard flujo de caja, exportar transacciones a Excel, despliegue en Vercel'
$versiones.Range("A3").NumberFormat = "@"
$versiones.Range("A3").Value = '1.1'
$versiones.Range("A3").Style = "Normal"
$versiones.Range("B3").Value = '27/02/2025'
$versiones.Range("C3").Value = 'Regla flujo despliegue (probar en local → confirmar → desplegar); hoja Versiones en bitácora'

Write-Output "edit applied"
